$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-12 Thursday" "2024-12-13 Friday"

Replace-Text "377×9=" "894×8="
Replace-Text "547×4=" "422×2="
Replace-Text "391×8=" "671×2="
Replace-Text "517×3=" "844×2="
Replace-Text "794×7=" "538×5="

Replace-Text "835×4=" "545×7="
Replace-Text "609×3=" "231×9="
Replace-Text "428×3=" "657×5="
Replace-Text "981×5=" "561×2="
Replace-Text "137×6=" "802×4="

Replace-Text "780×4=" "393×7="
Replace-Text "401×5=" "626×7="
Replace-Text "420×8=" "861×3="
Replace-Text "371×2=" "876×4="
Replace-Text "184×2=" "849×4="

Replace-Text "359×3=" "962×9="
Replace-Text "687×2=" "250×5="
Replace-Text "647×8=" "448×3="
Replace-Text "143×7=" "739×3="
Replace-Text "689×6=" "815×6="

Replace-Text "505×6=" "143×3="
Replace-Text "550×4=" "297×6="
Replace-Text "872×7=" "167×4="
Replace-Text "704×8=" "215×2="
Replace-Text "797×3=" "434×2="
